$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for the 01e0d4b6-...md row was regenerated,
# moving from 2016-09-07 15:50:13 to 2016-09-07 15:51:36. This datetime is
# shared by Overview!G4:G7 and de-de!H4:H7 (they all reference the same
# underlying value).
$wsOverview.Range("G4:G7").Value = "2016-09-07 15:51:36"
$wsDeDe.Range("H4:H7").Value = "2016-09-07 15:51:36"

# Priority column ("low" -> "ht") for rows 4-7 in both locale sheets.
$wsZhCn.Range("E4:E7").Value = "ht"
$wsDeDe.Range("E4:E7").Value = "ht"

# "Latest Handoff Datetime" for zh-cn rows 4-7 moved from
# 2016-09-07 15:49:52 to 2016-09-07 15:51:19.
$wsZhCn.Range("H4:H7").Value = "2016-09-07 15:51:19"
